# Adds the "Partitions With Given Difference" (CN/GFG) question row that was
# previously blank (row 5), matching the formatting already used by its
# sibling rows (left/top aligned, row 4 style; wrap-text, row 4 "E" style).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new question row (row 5) that was previously blank.
# C5 is written first so "Partitions With Given Difference" becomes shared
# string #16 and "CN/GFG" (written via A5/B5 next) becomes #17, matching
# the insertion order recorded in the target workbook.
$ws.Range("C5").Value = "Partitions With Given Difference"
$ws.Range("A5").Value = "CN/GFG"
$ws.Range("B5").Value = "CN/GFG"
$ws.Range("D5").Value = "Java"
$ws.Range("E5").Value = "DP(Recurrsion+Memonization+Tabulation+Space optimization)"

# A5:E5 share the same left/top, wrap-text look already used one row above
# (A4:E4) for this same "question" table.
$newRowRange = $ws.Range("A5:E5")
$newRowRange.HorizontalAlignment = -4131
$newRowRange.VerticalAlignment = -4160
$newRowRange.WrapText = $true

# Rows 4 and 5 grow to fit the new, longer two-line text.
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30

# The saved selection moves from C10 to D10.
$ws.Range("D10").Select()
